# Update scripts with new TPM values.
# The "ECs" sending-cluster row (original row 2) is removed, and the
# remaining rows are recalculated with the new TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for the "ECs" sending cluster (original row 2).
# This shifts the old row 3 (FAPs) up to row 2 and old row 4 (MuSCs) up to row 3,
# and Excel automatically drops the now-unused "ECs" shared string.
$ws.Rows("2").Delete()

# Row 2 (Sending cluster = FAPs, Ligand = Ccl17, Receptor = Ackr2, Target = FAPs)
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.4342616666666667
$ws.Range("H2").Value2 = 1.302785
$ws.Range("I2").Value2 = 0.8244337173731726
$ws.Range("J2").Value2 = 0.8244337173731725
$ws.Range("Q2").Value2 = 5.00149207419889
$ws.Range("R2").Value2 = 45.01342866779
$ws.Range("S2").Value2 = 0.8244337173731726
$ws.Range("T2").Value2 = 0.8244337173731725

# Row 3 (Sending cluster = MuSCs, Ligand = Ccl17, Receptor = Ackr2, Target = FAPs)
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.09247766666666667
$ws.Range("H3").Value2 = 0.277433
$ws.Range("I3").Value2 = 0.1755662826268274
$ws.Range("J3").Value2 = 0.1755662826268274
$ws.Range("Q3").Value2 = 1.065086680166889
$ws.Range("R3").Value2 = 9.585780121501999
$ws.Range("S3").Value2 = 0.1755662826268274
$ws.Range("T3").Value2 = 0.1755662826268274
